$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1277.3462
$ws.Range("I2").Value = 250.64285
$ws.Range("J2").Value = 2475.1667
$ws.Range("K2").Value = 250.64285
$ws.Range("L2").Value = 2475.1667
$ws.Range("M2").Value = -137.64285
$ws.Range("N2").Value = -2701.1667

$ws.Range("H9").Value = 700
$ws.Range("I9").Value = 468
$ws.Range("J9").Value = 1570
$ws.Range("K9").Value = 468
$ws.Range("L9").Value = 1570
$ws.Range("M9").Value = -299
$ws.Range("N9").Value = -1908

$ws.Range("H40").Value = 1943.8889
$ws.Range("I40").Value = 1870.8572
$ws.Range("J40").Value = 2199.5
$ws.Range("K40").Value = 1870.8572
$ws.Range("L40").Value = 2199.5
$ws.Range("M40").Value = -1695.8572
$ws.Range("N40").Value = -2549.5

$ws.Range("H64").Value = 4256.7144
$ws.Range("I64").Value = 3549.5
$ws.Range("K64").Value = 3549.5
$ws.Range("M64").Value = -3301.5

$ws.Range("H67").Value = 4256.7144
$ws.Range("I67").Value = 3549.5
$ws.Range("K67").Value = 3549.5
$ws.Range("M67").Value = -2691.5

$ws.Range("H88").Value = 4731.5713
$ws.Range("I88").Value = 569.6667
$ws.Range("K88").Value = 569.6667
$ws.Range("M88").Value = -163.6667

$ws.Range("H91").Value = 4731.5713
$ws.Range("I91").Value = 569.6667
$ws.Range("K91").Value = 569.6667
$ws.Range("M91").Value = 834.3333

$ws.Range("H92").Value = 190.93333
$ws.Range("I92").Value = 147.77777
$ws.Range("K92").Value = 147.77777
$ws.Range("M92").Value = 1100.22223

$ws.Range("H96").Value = 6689.625
$ws.Range("I96").Value = 8846.666999999999
$ws.Range("K96").Value = 26540.001
$ws.Range("M96").Value = -25167.001

$ws.Range("H98").Value = 3001
$ws.Range("J98").Value = 3901.7
$ws.Range("L98").Value = 3901.7
$ws.Range("N98").Value = -6897.7

$ws.Range("H112").Value = 3927.5
$ws.Range("I112").Value = 3600
$ws.Range("J112").Value = 3993
$ws.Range("K112").Value = 10800
$ws.Range("L112").Value = 11979
$ws.Range("M112").Value = -9692
$ws.Range("N112").Value = -14195

$ws.Range("H122").Value = 3001
$ws.Range("J122").Value = 3901.7
$ws.Range("L122").Value = 11705.1
$ws.Range("N122").Value = -16605.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 198
$ws.Range("I12").Value = 198
$ws.Range("K12").Value = 198
$ws.Range("M12").Value = -25

$ws.Range("H61").Value = 3520.3103
$ws.Range("I61").Value = 1670.3077
$ws.Range("J61").Value = 5023.4375
$ws.Range("K61").Value = 1670.3077
$ws.Range("L61").Value = 5023.4375
$ws.Range("M61").Value = -1458.3077
$ws.Range("N61").Value = -5447.4375

$ws.Range("H64").Value = 72500
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 72500
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 72500
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -72996

$ws.Range("H67").Value = 72500
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 72500
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 72500
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -74216

$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws.Range("H136").Value = 3520.3103
$ws.Range("I136").Value = 1670.3077
$ws.Range("J136").Value = 5023.4375
$ws.Range("K136").Value = 5010.9231
$ws.Range("L136").Value = 15070.3125
$ws.Range("M136").Value = -2460.9231
$ws.Range("N136").Value = -20170.3125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1832.25
$ws.Range("I94").Value = 1736.8572
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 1736.8572
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -1285.8572
$ws.Range("N94").Value = -3402

$ws.Range("H107").Value = 2999.5715
$ws.Range("I107").Value = 2999.5715
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2999.5715
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1079.5715
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = 84999
$ws.Range("J132").Value = 84999
$ws.Range("L132").Value = 84999
$ws.Range("M132").Value = -95119

$ws.Range("H134").Value = 1537.3334
$ws.Range("I134").Value = 1537.3334
$ws.Range("K134").Value = 4612.0002
$ws.Range("M134").Value = -2077.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 8305.111000000001
$ws.Range("J86").Value = 9999
$ws.Range("L86").Value = 9999
$ws.Range("N86").Value = -12245

$ws.Range("H89").Value = 8305.111000000001
$ws.Range("J89").Value = 9999
$ws.Range("L89").Value = 49995
$ws.Range("N89").Value = -61227

$ws.Range("H105").Value = 3222.0322
$ws.Range("I105").Value = 2436.2222
$ws.Range("K105").Value = 2436.2222
$ws.Range("M105").Value = -689.2222000000002

$ws.Range("H141").Value = 561324
$ws.Range("J141").Value = 708333.3
$ws.Range("L141").Value = 708333.3
$ws.Range("N141").Value = -718693.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 758.1429000000001
$ws.Range("J46").Value = 758.1429000000001
$ws.Range("L46").Value = 2274.4287
$ws.Range("N46").Value = -2456.4287

$ws.Range("H117").Value = 12990.25
$ws.Range("I117").Value = 425.2
$ws.Range("J117").Value = 33932
$ws.Range("K117").Value = 1275.6
$ws.Range("L117").Value = 101796
$ws.Range("M117").Value = 2166.4
$ws.Range("N117").Value = -108680

$ws.Range("H118").Value = 4014
$ws.Range("I118").Value = 4014
$ws.Range("K118").Value = 12042
$ws.Range("M118").Value = -10799

$ws.Range("H124").Value = 2783.25
$ws.Range("I124").Value = 1450
$ws.Range("J124").Value = 4116.5
$ws.Range("K124").Value = 4350
$ws.Range("L124").Value = 12349.5
$ws.Range("M124").Value = 560
$ws.Range("N124").Value = -22169.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1700
$ws.Range("I80").Value = 1700
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1700
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -702
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 1700
$ws.Range("I83").Value = 1700
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 8500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -3508
$ws.Range("N83").ClearContents()

$ws.Range("H122").Value = 608.8
$ws.Range("I122").Value = 598.75
$ws.Range("K122").Value = 1796.25
$ws.Range("M122").Value = 653.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2983.5557
$ws.Range("I46").Value = 2978.8572
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 2978.8572
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = -2790.8572
$ws.Range("N46").Value = -3376

$ws.Range("H100").Value = 3913
$ws.Range("I100").Value = 2992.6667
$ws.Range("J100").Value = 4833.3335
$ws.Range("K100").Value = 2992.6667
$ws.Range("L100").Value = 4833.3335
$ws.Range("M100").Value = -2451.6667
$ws.Range("N100").Value = -5915.3335

$ws.Range("H122").Value = 7874.75
$ws.Range("I122").Value = 9000
$ws.Range("J122").Value = 6544.909
$ws.Range("K122").Value = 27000
$ws.Range("L122").Value = 19634.727
$ws.Range("M122").Value = -24550
$ws.Range("N122").Value = -24534.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 90095
$ws.Range("I70").Value = 90095
$ws.Range("K70").Value = 90095
$ws.Range("M70").Value = -89780

$ws.Range("H73").Value = 90095
$ws.Range("I73").Value = 90095
$ws.Range("K73").Value = 90095
$ws.Range("M73").Value = -89003

$ws.Range("H100").Value = 5265167.5
$ws.Range("I100").Value = 6251524
$ws.Range("K100").Value = 12503048
$ws.Range("M100").Value = -12502507

$ws.Range("H122").Value = 1051.909
$ws.Range("I122").Value = 897.3333
$ws.Range("K122").Value = 2691.9999
$ws.Range("M122").Value = -241.9998999999998

$ws.Range("H126").Value = 5736.7646
$ws.Range("I126").Value = 4720
$ws.Range("J126").Value = 7600.8335
$ws.Range("K126").Value = 14160
$ws.Range("L126").Value = 22802.5005
$ws.Range("M126").Value = -11690
$ws.Range("N126").Value = -27742.5005
